$wb = $excel.ActiveWorkbook

# Each of these cells currently holds a bare cardinality fragment like "..1",
# "..*", "..0", or "..4" (the upper-bound half of a split cardinality value).
# The commit prefixes each of them with a single leading space.
$wb.Worksheets.Item("Coverage").Range("B5").Value = " ..1"
$wb.Worksheets.Item("Device").Range("B2").Value = " ..1"
$wb.Worksheets.Item("DocumentReference").Range("B7").Value = " ..1"
$wb.Worksheets.Item("Encounter").Range("B4").Value = " ..*"
$wb.Worksheets.Item("Location").Range("B2").Value = " ..*"
$wb.Worksheets.Item("Location").Range("C3").Value = " ..1"
$wb.Worksheets.Item("Observation").Range("AA43").Value = " ..0"
$wb.Worksheets.Item("Observation").Range("O51").Value = " ..0"
$wb.Worksheets.Item("Organization").Range("B3").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("B4").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("B5").Value = " ..4"
$wb.Worksheets.Item("Organization").Range("B6").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("B7").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("B8").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("B9").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("C10").Value = " ..1"
$wb.Worksheets.Item("Organization").Range("C13").Value = " ..1"
$wb.Worksheets.Item("Patient").Range("C6").Value = " ..1"
$wb.Worksheets.Item("Practitioner").Range("B2").Value = " ..1"
$wb.Worksheets.Item("Practitioner").Range("B3").Value = " ..1"
$wb.Worksheets.Item("Practitioner").Range("B4").Value = " ..4"
$wb.Worksheets.Item("Practitioner").Range("B5").Value = " ..1"
$wb.Worksheets.Item("Practitioner").Range("B6").Value = " ..1"
$wb.Worksheets.Item("Provenance").Range("B2").Value = " ..*"
